$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("B12").Value = 21
$ws.Range("B13").Value = 735000
$ws.Range("B14").Value = 2250000
$ws.Range("B33").Value = 1655000
$ws.Range("B35").Value = 1655000
